$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.682284333333333
$ws.Range("H2").Value = 14.046853
$ws.Range("I2").Value = 0.1747450949072225
$ws.Range("J2").Value = 0.1747450949072225
$ws.Range("Q2").Value = 0.1263920225325555
$ws.Range("R2").Value = 1.137528202793
$ws.Range("S2").Value = 0.1747450949072225
$ws.Range("T2").Value = 0.1747450949072225

$ws.Range("I3").Value = 0.6577126751045782
$ws.Range("J3").Value = 0.6577126751045781
$ws.Range("S3").Value = 0.6577126751045782
$ws.Range("T3").Value = 0.6577126751045781

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9374000000000001
$ws.Range("H4").Value = 2.8122
$ws.Range("I4").Value = 0.03498421716936109
$ws.Range("J4").Value = 0.03498421716936108
$ws.Range("Q4").Value = 0.02530386313333334
$ws.Range("R4").Value = 0.2277347682
$ws.Range("S4").Value = 0.03498421716936109
$ws.Range("T4").Value = 0.03498421716936108

$ws.Range("G5").Value = 3.238087
$ws.Range("H5").Value = 9.714261
$ws.Range("I5").Value = 0.1208469584182685
$ws.Range("J5").Value = 0.1208469584182685
$ws.Range("Q5").Value = 0.08740784111566667
$ws.Range("R5").Value = 0.7866705700410001
$ws.Range("S5").Value = 0.1208469584182685
$ws.Range("T5").Value = 0.1208469584182685

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.313797
$ws.Range("H6").Value = 0.941391
$ws.Range("I6").Value = 0.01171105440056966
$ws.Range("J6").Value = 0.01171105440056966
$ws.Range("Q6").Value = 0.008470531618999999
$ws.Range("R6").Value = 0.076234784571
$ws.Range("S6").Value = 0.01171105440056966
$ws.Range("T6").Value = 0.01171105440056966
